$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.243.71"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "2.580.91"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'590.12"
$ws.Range("E5").Value = "  -2.82%  "
$ws.Range("D6").Value = "'150.01"
$ws.Range("E6").Value = "  +2.09%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "'5.69"
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "'0.150"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("D13").Value = "'27.42"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "3.043.84"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "63.021.23"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "'0.0000155"
$ws.Range("E16").Value = "  +5.84%  "
$ws.Range("D17").Value = "2.579.36"
$ws.Range("E17").Value = "  -1.82%  "
$ws.Range("D18").Value = "'12.23"
$ws.Range("E18").Value = "  +4.11%  "
$ws.Range("D19").Value = "'4.76"
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("D20").Value = "'344.29"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").Value = "'6.86"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "'9.25"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").Value = "'563.70"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'8.00"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.162"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "'2.03"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "0.0₃0844"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "'1.75"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").Value = "'5.21"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").Value = "'166.45"
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("D36").Value = "'0.411"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'19.46"
$ws.Range("E38").Value = "  +1.83%  "
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'166.73"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("D42").Value = "'39.53"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("E43").Value = "  +3.99%  "
$ws.Range("D44").Value = "'22.62"
$ws.Range("E44").Value = "  +3.77%  "
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.06"
$ws.Range("E46").Value = "  +3.09%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.628"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  +3.02%  "
$ws.Range("D49").Value = "'0.0961"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").Value = "'18.98"
$ws.Range("E50").Value = "  +1.14%  "
$ws.Range("E51").Value = "  +18.14%  "
